$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the old cells that are no longer used
$ws.Range("C2").ClearContents()
$ws.Range("B4").ClearContents()

# Write the new values into column A
$ws.Range("A2").Value = "ABCD"
$ws.Range("A3").Value = "EFG"
$ws.Range("A4").Value = "HIJ"
$ws.Range("A5").Value = "QRS"

# Update the selection to match the target state
$ws.Range("A5:A7").Select()
